# Adds a new week of attendance data (and fills a couple of catch-up cells
# for a prior week) across the Astronauta, Senador and Ninja sheets, then
# leaves Senador as the active/selected sheet (it was Ninja before).

$wb = $excel.ActiveWorkbook

$wsAstronauta = $wb.Worksheets.Item("Astronauta")
$wsSenador    = $wb.Worksheets.Item("Senador")
$wsMago       = $wb.Worksheets.Item("Mago")
$wsNinja      = $wb.Worksheets.Item("Ninja")

# ---------------------------------------------------------------------------
# Astronauta: new week column G (date 2023-11-24 / serial 45229).
# G1 already existed as an empty, rotated placeholder cell, so copy the
# date-formatted style from F1 onto it before writing the new date value.
# ---------------------------------------------------------------------------
$wsAstronauta.Range("F1").Copy()
$wsAstronauta.Range("G1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$wsAstronauta.Range("G1").Value = 45229

$wsAstronauta.Range("G5").Value = 1
$wsAstronauta.Range("G9").Value = 1
$wsAstronauta.Range("F17").Value = 1
$wsAstronauta.Range("F22").Value = 1
$wsAstronauta.Range("G23").Value = 1
$wsAstronauta.Range("G25").Value = 1

$wsAstronauta.Range("G2").Select()

# ---------------------------------------------------------------------------
# Senador: new week column G (date 2023-11-24 / serial 45229). G1 never
# existed on this sheet, so give it a plain (non-rotated) date format.
# ---------------------------------------------------------------------------
$wsSenador.Range("G1").Value = 45229
$wsSenador.Range("G1").NumberFormat = "mm-dd-yy"

$wsSenador.Range("F3").Value = 0
$wsSenador.Range("G3").Value = 0
$wsSenador.Range("F12").Value = 0
$wsSenador.Range("F15").Value = 0
$wsSenador.Range("G15").Value = 0
$wsSenador.Range("G23").Value = 0
$wsSenador.Range("G28").Value = 0

# ---------------------------------------------------------------------------
# Ninja: new week column J. J1 only gets the new date-number-format style
# (copied from Senador's freshly-styled G1 so both sheets share one style
# entry instead of minting a duplicate) but no value is entered yet.
# ---------------------------------------------------------------------------
$wsSenador.Range("G1").Copy()
$wsNinja.Range("J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsNinja.Range("J2").Value = 0
$wsNinja.Range("J5").Value = 1
$wsNinja.Range("J9").Value = 1
$wsNinja.Range("H13").Value = 0
$wsNinja.Range("H15").Value = 0
$wsNinja.Range("J17").Value = 1
$wsNinja.Range("I19").Value = 0
$wsNinja.Range("J22").Value = 1
$wsNinja.Range("J23").Value = 1
$wsNinja.Range("J25").Value = 1
$wsNinja.Range("J26").Value = 0

$wsNinja.Range("J18").Select()

# ---------------------------------------------------------------------------
# Mago is untouched by this edit. Make Senador the active sheet (it replaces
# Ninja as the workbook's selected tab) and restore its own selection.
# ---------------------------------------------------------------------------
$wsSenador.Activate()
$wsSenador.Range("H4").Select()
